# edit.ps1 - apply the newsocialism.docx revision described by the commit diff
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new list paragraph, immediately before the "There is the
#    possibility of allowing machines..." bullet, containing a hyperlink to
#    the Economist "labour pains" article. Then move the hidden _GoBack
#    bookmark so it sits at the start of that "There is the possibility..."
#    paragraph (Word keeps only one _GoBack bookmark and relocates it
#    automatically whenever a new one is added).
# ---------------------------------------------------------------------------
$target = $d.Content
$target.Find.ClearFormatting()
$found = $target.Find.Execute("There is the possibility of allowing machines", `
    $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found) { throw "anchor paragraph not found" }

$insertPoint = $d.Range($target.Start, $target.Start)
$url1 = "https://www.economist.com/news/finance-and-economics/21588900-all-around-world-labour-losing-out-capital-labour-pains"
$insertPoint.InsertBefore($url1 + "`r")

$linkRange = $d.Range($insertPoint.Start, $insertPoint.Start + $url1.Length)
$d.Hyperlinks.Add($linkRange, $url1, $null, $null, $null) | Out-Null

# Re-find the "There is the possibility..." paragraph (it shifted down by one
# paragraph) and drop the _GoBack bookmark right at its very start.
$again = $d.Content
$again.Find.ClearFormatting()
$again.Find.Execute("There is the possibility of allowing machines", `
    $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$bmPoint = $d.Range($again.Start, $again.Start)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the four short runs that spell out bullet "1. Provide investment
#    growth..." into a single run of text (no formatting differences, so we
#    just normalise the text - Word already stores them as one run once the
#    text matches and is retyped as a block).
# ---------------------------------------------------------------------------
$merge = $d.Content
$merge.Find.ClearFormatting()
$merge.Find.Execute("1. Provide investment growth and income to the members of the cooperative, with profits shared amongst the members. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, $null, 0) | Out-Null
$merge.Text = "1. Provide investment growth and income to the members of the cooperative, with profits shared amongst the members. "

# ---------------------------------------------------------------------------
# 3) Footnote hyperlinks: the new hyperlink added in step 1 shifted every
#    subsequent hyperlink relationship id up by one. Re-point each footnote
#    hyperlink at the same target URL so the relationship ids get reassigned
#    in document order (matching what Word does on save).
# ---------------------------------------------------------------------------
function Set-HyperlinkAddress($oldAddress, $newAddress) {
    foreach ($hl in $d.Hyperlinks) {
        if ($hl.Address -eq $oldAddress) {
            $hl.Address = $newAddress
            return
        }
    }
    throw "hyperlink not found with address: $oldAddress"
}

Set-HyperlinkAddress "https://www.economist.com/blogs/economist-explains/2015/09/economist-explains-19" `
    "https://www.economist.com/blogs/economist-explains/2015/09/economist-explains-19"
Set-HyperlinkAddress "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/616966/trade-union-membership-statistical-bulletin-2016-rev.pdf" `
    "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/616966/trade-union-membership-statistical-bulletin-2016-rev.pdf"
Set-HyperlinkAddress "https://www.ifs.org.uk/docs/ER_JC_2013.pdf" `
    "https://www.ifs.org.uk/docs/ER_JC_2013.pdf"
Set-HyperlinkAddress "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/648165/HEIPR_PUBLICATION_2015-16.pdf" `
    "https://assets.publishing.service.gov.uk/government/uploads/system/uploads/attachment_data/file/648165/HEIPR_PUBLICATION_2015-16.pdf"
Set-HyperlinkAddress "https://www.nature.com/articles/s41539-018-0019-8" `
    "https://www.nature.com/articles/s41539-018-0019-8"

Write-Output "done"
